$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")
$ws.Range("B11").Value = "1"
